# New crime data collected - weekly 115th Precinct CompStat update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text updates (shared-string rich-text runs in merged header cells)
# ---------------------------------------------------------------------------
# A8: "Volume 32   Number  9" -> "Volume 32   Number  10" (last run "9" -> "10")
$ws.Range("A8").Characters(21, 1).Text = "10"

# C9: "Report Covering the Week  2/24/2025  Through  3/2/2025"
#  -> "Report Covering the Week  3/3/2025  Through  3/9/2025"
# Replace the right-most date first so the left-most replacement (which
# changes length) doesn't shift the already-computed character offsets.
$ws.Range("C9").Characters(47, 8).Text = "3/9/2025"
$ws.Range("C9").Characters(27, 9).Text = "3/3/2025"

# ---------------------------------------------------------------------------
# 2) Plain numeric value updates (value only, style/format unchanged)
# ---------------------------------------------------------------------------
$numericChanges = @{
    "N14" = -50
    "C16" = 5
    "D16" = 2
    "E16" = 150
    "F16" = 21
    "G16" = 24
    "H16" = -12.5
    "I16" = 42
    "J16" = 66
    "K16" = -36.363636363636
    "L16" = -8.695652173913
    "M16" = -27.586206896551
    "N16" = -82.5
    "C17" = 7
    "D17" = 6
    "E17" = 16.666666666666
    "F17" = 34
    "G17" = 38
    "H17" = -10.526315789473
    "I17" = 85
    "J17" = 87
    "K17" = -2.298850574712
    "L17" = -5.555555555555
    "M17" = 73.469387755102
    "N17" = 34.920634920634
    "C18" = 2
    "D18" = 6
    "E18" = -66.666666666666
    "F18" = 12
    "G18" = 18
    "H18" = -33.333333333333
    "I18" = 30
    "J18" = 44
    "K18" = -31.818181818181
    "L18" = 57.894736842105
    "M18" = -43.396226415094
    "N18" = -93.814432989690
    "C19" = 11
    "D19" = 19
    "E19" = -42.105263157894
    "F19" = 53
    "G19" = 58
    "H19" = -8.620689655172
    "I19" = 132
    "J19" = 184
    "K19" = -28.260869565217
    "L19" = -11.409395973154
    "M19" = 55.294117647058
    "N19" = -43.589743589743
    "C20" = 1
    "D20" = 6
    "E20" = -83.333333333333
    "F20" = 13
    "H20" = 0
    "I20" = 31
    "J20" = 44
    "K20" = -29.545454545454
    "L20" = -46.551724137931
    "M20" = -26.190476190476
    "N20" = -92.420537897310
    "C21" = 26
    "D21" = 39
    "E21" = -33.333333333333
    "F21" = 133
    "H21" = -12.5
    "I21" = 328
    "J21" = 429
    "K21" = -23.543123543123
    "L21" = -10.382513661202
    "M21" = 13.103448275862
    "N21" = -77.269577269577
    "D22" = 3
    "E22" = -66.666666666666
    "G22" = 4
    "H22" = 0
    "I22" = 10
    "J22" = 13
    "K22" = -23.076923076923
    "L22" = -37.5
    "M22" = 66.666666666666
    "C24" = 19
    "D24" = 42
    "E24" = -54.761904761904
    "F24" = 98
    "G24" = 215
    "H24" = -54.418604651162
    "I24" = 253
    "J24" = 476
    "K24" = -46.848739495798
    "L24" = -41.162790697674
    "M24" = 40.555555555555
    "C25" = 10
    "D25" = 21
    "E25" = -52.380952380952
    "F25" = 41
    "G25" = 128
    "H25" = -67.96875
    "I25" = 117
    "J25" = 289
    "K25" = -59.515570934256
    "L25" = -50.632911392405
    "C26" = 31
    "D26" = 32
    "E26" = -3.125
    "F26" = 75
    "G26" = 101
    "H26" = -25.742574257425
    "I26" = 193
    "J26" = 206
    "K26" = -6.310679611650
    "L26" = 27.814569536423
    "M26" = 13.529411764705
    "G27" = 3
    "H27" = -66.666666666666
    "J27" = 11
    "K27" = -27.272727272727
    "C28" = 2
    "D28" = 3
    "E28" = -33.333333333333
    "F28" = 10
    "G28" = 9
    "H28" = 11.111111111111
    "I28" = 17
    "J28" = 25
    "K28" = -32
    "L28" = 0
    "J31" = 4
}

foreach ($ref in $numericChanges.Keys) {
    $ws.Range($ref).Value = $numericChanges[$ref]
}

# ---------------------------------------------------------------------------
# 3) Cells that flip between the "N/A" text marker and a real number
#    (style also changes, e.g. General <-> #,##0 / #,##0.0). We copy the
#    formatting+content from a donor cell that already carries the target
#    style in the same row, then overwrite the value where needed, so both
#    the stored style index and shared-string usage line up with the target.
# ---------------------------------------------------------------------------

# C22: "N/A" text -> number 1 (copy numeric style from D22)
$ws.Range("D22").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1

# D27: "N/A" text -> number 1 (copy numeric style from F27)
$ws.Range("F27").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
# E27: "***.*" text -> number -100 (copy percent style from H27)
$ws.Range("H27").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100

# D31: "N/A" text -> number 1 (copy numeric style from G31)
$ws.Range("G31").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
# E31: "***.*" text -> number -100 (copy percent style from H31)
$ws.Range("H31").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100

# C33: number 1 -> "N/A" text (copy text style+value from D33)
$ws.Range("D33").Copy($ws.Range("C33"))
# G33: number 1 -> "N/A" text (copy text style+value from D33)
$ws.Range("D33").Copy($ws.Range("G33"))
# H33: number 0 -> "***.*" text (copy text style+value from E33)
$ws.Range("E33").Copy($ws.Range("H33"))
